$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich-text cells) ---
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# --- Crime-statistics grid updates (rows 14-30) ---
$ws.Range("F14").Value = 1
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -72.413793103448
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 6
$ws.Range("I15").Value = 33
$ws.Range("K15").Value = 13.793103448275
$ws.Range("L15").Value = 26.923076923076
$ws.Range("M15").Value = 43.478260869565
$ws.Range("N15").Value = -48.4375
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = 5.263157894736
$ws.Range("I16").Value = 191
$ws.Range("J16").Value = 203
$ws.Range("K16").Value = -5.911330049261
$ws.Range("L16").Value = -28.195488721804
$ws.Range("M16").Value = -45.892351274787
$ws.Range("N16").Value = -89.054441260745
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -46.666666666666
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = -6.521739130434
$ws.Range("I17").Value = 516
$ws.Range("J17").Value = 549
$ws.Range("K17").Value = -6.010928961748
$ws.Range("L17").Value = -8.021390374331
$ws.Range("M17").Value = 9.322033898305
$ws.Range("N17").Value = -42.281879194630
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 109
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = 5.825242718446
$ws.Range("L18").Value = -26.351351351351
$ws.Range("M18").Value = -51.982378854625
$ws.Range("N18").Value = -82.305194805194
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 4.166666666666
$ws.Range("I19").Value = 278
$ws.Range("J19").Value = 259
$ws.Range("K19").Value = 7.335907335907
$ws.Range("L19").Value = -9.150326797385
$ws.Range("M19").Value = -17.261904761904
$ws.Range("N19").Value = -62.330623306233
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = -0.854700854700
$ws.Range("L20").Value = -22.666666666666
$ws.Range("M20").Value = -0.854700854700
$ws.Range("N20").Value = -77.299412915851
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -19.354838709677
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -0.892857142857
$ws.Range("I21").Value = 1259
$ws.Range("J21").Value = 1266
$ws.Range("K21").Value = -0.552922590837
$ws.Range("L21").Value = -14.528173794976
$ws.Range("M21").Value = -18.721755971594
$ws.Range("N21").Value = -72.784262862083
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 36.842105263157
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -46.428571428571
$ws.Range("I23").Value = 266
$ws.Range("J23").Value = 277
$ws.Range("K23").Value = -3.971119133574
$ws.Range("L23").Value = -10.135135135135
$ws.Range("M23").Value = 36.410256410256
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 26.666666666666
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 67.307692307692
$ws.Range("I24").Value = 877
$ws.Range("J24").Value = 833
$ws.Range("K24").Value = 5.282112845138
$ws.Range("L24").Value = 0.114155251141
$ws.Range("M24").Value = 16.777629826897
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = 133.333333333333
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 219
$ws.Range("J25").Value = 222
$ws.Range("K25").Value = -1.351351351351
$ws.Range("L25").Value = -19.780219780219
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -9.523809523809
$ws.Range("F26").Value = 80
$ws.Range("G26").Value = 85
$ws.Range("H26").Value = -5.882352941176
$ws.Range("I26").Value = 683
$ws.Range("J26").Value = 765
$ws.Range("K26").Value = -10.718954248366
$ws.Range("L26").Value = -6.566347469220
$ws.Range("M26").Value = -33.495618305744
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 44
$ws.Range("K27").Value = 22.222222222222
$ws.Range("L27").Value = 7.317073170731
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -42.857142857142
$ws.Range("I28").Value = 54
$ws.Range("J28").Value = 67
$ws.Range("K28").Value = -19.402985074626
$ws.Range("L28").Value = -18.181818181818
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80
$ws.Range("J29").Value = 37
$ws.Range("K29").Value = 32.432432432432
$ws.Range("L29").Value = -26.865671641791
$ws.Range("M29").Value = -46.739130434782
$ws.Range("N29").Value = -80.784313725490
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -66.666666666666
$ws.Range("J30").Value = 33
$ws.Range("K30").Value = 21.212121212121
$ws.Range("L30").Value = -23.076923076923
$ws.Range("M30").Value = -45.945945945945
$ws.Range("N30").Value = -83.050847457627

"edit applied: " + 160 + " cell(s) updated"
